$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 415
$ws.Range("A415").Value = "E M E F MARIA DE LOURDES CASADINI DA SILVA"
$ws.Range("B415").Value = "RUA ANTONIO ALVES DE CARVALHO, 56 ESCOLA CASADINI. EXPANSAO. 68560-000 Santana do Araguaia - PA."
$ws.Range("C415").NumberFormat = "@"
$ws.Range("C415").Value = "94."
$ws.Range("D415").NumberFormat = "@"
$ws.Range("D415").Value = "`n15524655"
$ws.Range("E415").Value = "`nUrbana"
$ws.Range("F415").Value = "`nMunicipal"
$ws.Range("G415").Value = "`n`n                        Ensino Infantil, Ensino Fundamental, Anos Iniciais, Anos Finais"

# Row 416
$ws.Range("A416").Value = "E M E F VILA NOVA"
$ws.Range("B416").Value = "VILA NOVA MARO - RIO ARAPIUNS, 68115-000 Santarém - PA."
$ws.Range("C416").NumberFormat = "@"
$ws.Range("C416").Value = "Informação indisponível"
$ws.Range("D416").NumberFormat = "@"
$ws.Range("D416").Value = "`n15156770"
$ws.Range("E416").Value = "`nRural"
$ws.Range("F416").Value = "`nMunicipal"
$ws.Range("G416").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 417
$ws.Range("A417").Value = "EEEM ALVARO ADOLFO DA SILVEIRA"
$ws.Range("B417").Value = "AV MARECHAL RONDON, SN SANTA CLARA. 68005-120 Santarém - PA."
$ws.Range("C417").NumberFormat = "@"
$ws.Range("C417").Value = "(93) 3522-2329"
$ws.Range("D417").NumberFormat = "@"
$ws.Range("D417").Value = "`n15011372"
$ws.Range("E417").Value = "`nUrbana"
$ws.Range("F417").Value = "`nEstadual"
$ws.Range("G417").Value = "`n`n                        Ensino Médio"

# Row 418
$ws.Range("A418").Value = "E M E F NSRA DAS GRACAS"
$ws.Range("B418").Value = "COMUNIDADE DE AMARI, ARAPIUNS. 68115-000 Santarém - PA."
$ws.Range("C418").NumberFormat = "@"
$ws.Range("C418").Value = ".."
$ws.Range("D418").NumberFormat = "@"
$ws.Range("D418").Value = "`n15013260"
$ws.Range("E418").Value = "`nRural"
$ws.Range("F418").Value = "`nMunicipal"
$ws.Range("G418").Value = "`n`n                        Ensino Fundamental"

# Row 419
$ws.Range("A419").Value = "E M E F SANTISSIMA TRINDADE"
$ws.Range("B419").Value = "SANTISSIMA TRINDADE- VILA DO ARITAPERA, VARZEA. 68124-000 Santarém - PA."
$ws.Range("C419").NumberFormat = "@"
$ws.Range("C419").Value = ".."
$ws.Range("D419").NumberFormat = "@"
$ws.Range("D419").Value = "`n15140849"
$ws.Range("E419").Value = "`nRural"
$ws.Range("F419").Value = "`nMunicipal"
$ws.Range("G419").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 420
$ws.Range("A420").Value = "E M E I E F PROF OLINDO LUIZ DO CARMO NEVES"
$ws.Range("B420").Value = "RUA ITUQUI, S/N IGREJA DO AMPARO. AMPARO. 68035-670 Santarém - PA."
$ws.Range("C420").NumberFormat = "@"
$ws.Range("C420").Value = "(93) 99121-8227"
$ws.Range("D420").NumberFormat = "@"
$ws.Range("D420").Value = "`n15013626"
$ws.Range("E420").Value = "`nUrbana"
$ws.Range("F420").Value = "`nMunicipal"
$ws.Range("G420").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 421
$ws.Range("A421").Value = "E M E F VINTE DE JULHO"
$ws.Range("B421").Value = "COMUNIDADE CORREIO DO TAPARA, VARZEA. 68124-000 Santarém - PA."
$ws.Range("C421").NumberFormat = "@"
$ws.Range("C421").Value = "Informação indisponível"
$ws.Range("D421").NumberFormat = "@"
$ws.Range("D421").Value = "`n15147991"
$ws.Range("E421").Value = "`nRural"
$ws.Range("F421").Value = "`nMunicipal"
$ws.Range("G421").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 422
$ws.Range("A422").Value = "E M E F CEL MARIO FERNANDES IMBIRIBA"
$ws.Range("B422").Value = "RUA CASTELO BRANCO, INTERVENTORIA. 68015-260 Santarém - PA."
$ws.Range("C422").NumberFormat = "@"
$ws.Range("C422").Value = "(93) 3523-4313"
$ws.Range("D422").NumberFormat = "@"
$ws.Range("D422").Value = "`n15012115"
$ws.Range("E422").Value = "`nUrbana"
$ws.Range("F422").Value = "`nMunicipal"
$ws.Range("G422").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 423
$ws.Range("A423").Value = "E M E I E F NOVA ESPERANCA"
$ws.Range("B423").Value = "COMUNIDADE NOVA ESPERANCA DO ITUQUI, S/N ZONA RURAL. PLANALTO. 68128-000 Santarém - PA."
$ws.Range("C423").NumberFormat = "@"
$ws.Range("C423").Value = "Informação indisponível"
$ws.Range("D423").NumberFormat = "@"
$ws.Range("D423").Value = "`n15540650"
$ws.Range("E423").Value = "`nRural"
$ws.Range("F423").Value = "`nMunicipal"
$ws.Range("G423").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 424
$ws.Range("A424").Value = "E M E F JOSE DE MELO FILHO"
$ws.Range("B424").Value = "VILA DE AMORIM, RIO TAPAJOS. ZONA RURAL. 68115-000 Santarém - PA."
$ws.Range("C424").NumberFormat = "@"
$ws.Range("C424").Value = "(93) 3584-4125"
$ws.Range("D424").NumberFormat = "@"
$ws.Range("D424").Value = "`n15015955"
$ws.Range("E424").Value = "`nRural"
$ws.Range("F424").Value = "`nMunicipal"
$ws.Range("G424").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 425
$ws.Range("A425").Value = "E M E F RAIMUNDA DE LIRA MAIA"
$ws.Range("B425").Value = "TRAVESSA B, S/N ELCIONE BARBALHO. 68040-050 Santarém - PA."
$ws.Range("C425").NumberFormat = "@"
$ws.Range("C425").Value = "(93) 99182-2030"
$ws.Range("D425").NumberFormat = "@"
$ws.Range("D425").Value = "`n15567044"
$ws.Range("E425").Value = "`nUrbana"
$ws.Range("F425").Value = "`nMunicipal"
$ws.Range("G425").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 426
$ws.Range("A426").Value = "EEEF RICHARD HENNINGTON"
$ws.Range("B426").Value = "TRAVESSA XINGU, 997 ENTRE AVENIDA PALHAO. DIAMANTINO. 68020-140 Santarém - PA."
$ws.Range("C426").NumberFormat = "@"
$ws.Range("C426").Value = "(93) 3524-3435"
$ws.Range("D426").NumberFormat = "@"
$ws.Range("D426").Value = "`n15011712"
$ws.Range("E426").Value = "`nUrbana"
$ws.Range("F426").Value = "`nEstadual"
$ws.Range("G426").Value = "`n`n                        Ensino Fundamental"

# Row 427
$ws.Range("A427").Value = "E M E F SAO TOME"
$ws.Range("B427").Value = "COMUNIDADE SAO PEDRO, PLANALTO. 68123-000 Santarém - PA."
$ws.Range("C427").NumberFormat = "@"
$ws.Range("C427").Value = "(93) 3596-2063"
$ws.Range("D427").NumberFormat = "@"
$ws.Range("D427").Value = "`n15589404"
$ws.Range("E427").Value = "`nRural"
$ws.Range("F427").Value = "`nMunicipal"
$ws.Range("G427").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 428
$ws.Range("A428").Value = "E M E F NOSSA SRA DE FATIMA"
$ws.Range("B428").Value = "ALDEIA NOVA VISTA- ARAPIUNS, ARAPIUNS. 68115-000 Santarém - PA."
$ws.Range("C428").NumberFormat = "@"
$ws.Range("C428").Value = "Informação indisponível"
$ws.Range("D428").NumberFormat = "@"
$ws.Range("D428").Value = "`n15013049"
$ws.Range("E428").Value = "`nRural"
$ws.Range("F428").Value = "`nMunicipal"
$ws.Range("G428").Value = "`n`n                        Ensino Infantil, Ensino Fundamental, Anos Finais"

# Row 429
$ws.Range("A429").Value = "E M E F JAYME BARCESSAT"
$ws.Range("B429").Value = "HIDRELETRICA DE CURUA-UNA, PA 370 KM74, PLANALTO. 68010-000 Santarém - PA."
$ws.Range("C429").NumberFormat = "@"
$ws.Range("C429").Value = "93."
$ws.Range("D429").NumberFormat = "@"
$ws.Range("D429").Value = "`n15011232"
$ws.Range("E429").Value = "`nRural"
$ws.Range("F429").Value = "`nMunicipal"
$ws.Range("G429").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"

# Row 430
$ws.Range("A430").Value = "EEEM JULIA PASSARINHO - ANEXO I"
$ws.Range("B430").Value = "AVENIDA COSTA E SILVA, SN PROX. POSTO DE SAUD. MARARU. 68050-070 Santarém - PA."
$ws.Range("C430").NumberFormat = "@"
$ws.Range("C430").Value = "(93) 3523-5994"
$ws.Range("D430").NumberFormat = "@"
$ws.Range("D430").Value = "`n15170080"
$ws.Range("E430").Value = "`nUrbana"
$ws.Range("F430").Value = "`nEstadual"
$ws.Range("G430").Value = "`n`n                        Ensino Médio"

# Row 431
$ws.Range("A431").Value = "E M E F STA TEREZINHA"
$ws.Range("B431").Value = "BOCA DE CIMA DO ARITAPERA, VARZEA. 68124-000 Santarém - PA."
$ws.Range("C431").NumberFormat = "@"
$ws.Range("C431").Value = "(93) 99163-8872"
$ws.Range("D431").NumberFormat = "@"
$ws.Range("D431").Value = "`n15015084"
$ws.Range("E431").Value = "`nRural"
$ws.Range("F431").Value = "`nMunicipal"
$ws.Range("G431").Value = "`n`n                        Ensino Infantil, Ensino Fundamental"
